# Adding a new place near the Cascade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Cozy Apartments Near Cascade In DownTown"
$ws.Cells.Item($row, 3).Value = "Apartment"
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 40.191701778676297
$ws.Cells.Item($row, 6).Value = 44.517007512745899
$ws.Cells.Item($row, 7).Value = "Kentron"

$ws.Columns.Item(2).ColumnWidth = 41.6

$ws.Range("D15").Select()
